# Normalize the "Recorded By" column (G): when the list of recorders
# starts with "System" (case-insensitive), move it so that the other
# recorder(s) come first and "System" is listed last.
#
# Example:
#   "System, dnasr281@gmail.com"          -> "dnasr281@gmail.com, System"
#   "System, backup@backdoor.com, system" -> "system, backup@backdoor.com, System"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # column G = "Recorded By"
    $raw = $cell.Value2

    if ($raw -eq $null) {
        continue
    }

    $text = [string]$raw
    if ($text -eq "") {
        continue
    }

    $parts = $text -split ", "
    if ($parts.Count -le 1) {
        continue
    }

    $first = $parts[0].Trim()
    if ($first.ToLower() -ne "system") {
        continue
    }

    $reversed = $parts[($parts.Count - 1)..0]
    $newText = [string]::Join(", ", $reversed)

    if ($newText.CompareTo($text) -ne 0) {
        $cell.Value = $newText
    }
}
